# Daily attendance processing - rotate "Recorded By" (column G) values so that
# the last comma-separated entry moves to the front of the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Text

    if ([string]::IsNullOrEmpty($val)) {
        continue
    }

    $parts = $val -split ', '

    if ($parts.Length -gt 1) {
        $lastPart = $parts[$parts.Length - 1]
        $rest = $parts[0..($parts.Length - 2)]
        $newParts = @($lastPart) + $rest
        $newVal = $newParts -join ', '
        $cell.Value = $newVal
    }
}
